$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ECs" target-cluster row (original row 2). This shifts the
# remaining rows (FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac) up by one,
# and also drops the now-unused "ECs" shared string.
$ws.Rows("2:2").Delete()

# Update the recalculated (new TPM-based) numeric columns K:T for the
# four remaining data rows (now rows 2-5).

# Row 2 -> FAPs
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.257301333333333
$ws.Range("N2").Value = 3.771904
$ws.Range("O2").Value = 0.5503138561234041
$ws.Range("P2").Value = 0.5503138561234042
$ws.Range("Q2").Value = 0.1858723044124444
$ws.Range("R2").Value = 1.672850739712
$ws.Range("S2").Value = 0.5503138561234041
$ws.Range("T2").Value = 0.5503138561234042

# Row 3 -> Inflammatory-Mac
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.364727
$ws.Range("N3").Value = 1.094181
$ws.Range("O3").Value = 0.159638995426968
$ws.Range("P3").Value = 0.159638995426968
$ws.Range("Q3").Value = 0.05391917289366666
$ws.Range("R3").Value = 0.4852725560429999
$ws.Range("S3").Value = 0.159638995426968
$ws.Range("T3").Value = 0.159638995426968

# Row 4 -> MuSCs (K and L unchanged)
$ws.Range("M4").Value = 0.3637273333333333
$ws.Range("N4").Value = 1.091182
$ws.Range("O4").Value = 0.1592014468428805
$ws.Range("P4").Value = 0.1592014468428805
$ws.Range("Q4").Value = 0.05377138783844444
$ws.Range("R4").Value = 0.483942490546
$ws.Range("S4").Value = 0.1592014468428805
$ws.Range("T4").Value = 0.1592014468428805

# Row 5 -> Resolving-Mac
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.298943
$ws.Range("N5").Value = 0.896829
$ws.Range("O5").Value = 0.1308457016067473
$ws.Range("P5").Value = 0.1308457016067473
$ws.Range("Q5").Value = 0.04419403910966666
$ws.Range("R5").Value = 0.397746351987
$ws.Range("S5").Value = 0.1308457016067473
$ws.Range("T5").Value = 0.1308457016067473
